$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Bibi Cell Mundi
$ws.Range("AA2").Value = 15187.76
$ws.Range("AG2").Value = 290748.54

# Row 3 - Bibi Cell Vieiralves
$ws.Range("AA3").Value = 5175
$ws.Range("AG3").Value = 173119.3

# Row 4 - Bibi Cell Ponta Negra
$ws.Range("Z4").Value = 1556
$ws.Range("AA4").Value = 5059.01
$ws.Range("AG4").Value = 87970.33

# Row 5 - Bibi Cell Manauara
$ws.Range("Z5").Value = 2268.9
$ws.Range("AA5").Value = 1972
$ws.Range("AG5").Value = 79827.48

# Row 6 - total
$ws.Range("Z6").Value = 3824.9
$ws.Range("AA6").Value = 27393.77
$ws.Range("AG6").Value = 631665.65
